$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the "ВЕЙП" and "ВЕЙПЕР" rows (rows 16 and 17), shifting the rows
# below them up. Column A keeps its original running numbering (which
# continues counting the deleted rows), so after the delete we restore it.
$ws.Rows("16:17").Delete()

# Restore column-A running numbers for the rows that shifted up, so they
# continue the sequence as if rows 16/17 had simply been skipped (18..25)
for ($i = 16; $i -le 23; $i++) {
    $ws.Cells.Item($i, 1).Value = $i + 2
}

$ws.Range("D20").Select() | Out-Null
